$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "51.490.77"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "3.038.67"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("E4").Value = "  +0.10%  "

Set-TextValue $ws.Range("D5") "384.32"
$ws.Range("E5").Value = "  +1.26%  "

Set-TextValue $ws.Range("D6") "102.63"
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue $ws.Range("D9") "0.583"
$ws.Range("E9").Value = "  -1.12%  "

Set-TextValue $ws.Range("D10") "36.80"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "3.526.61"
$ws.Range("E13").Value = "  +3.11%  "

Set-TextValue $ws.Range("D14") "18.57"
$ws.Range("E14").Value = "  +1.77%  "

Set-TextValue $ws.Range("D15") "7.72"
$ws.Range("E15").Value = "  -1.01%  "

$ws.Range("D16").Value = "3.042.71"
$ws.Range("E16").Value = "  +3.26%  "

Set-TextValue $ws.Range("D17") "0.970"
$ws.Range("E17").Value = "  -2.69%  "

Set-TextValue $ws.Range("D18") "10.46"
$ws.Range("E18").Value = "  -6.10%  "

$ws.Range("D19").Value = "51.549.90"
$ws.Range("E19").Value = "  +0.81%  "

Set-TextValue $ws.Range("D20") "3.13"
$ws.Range("E20").Value = "  -1.82%  "

Set-TextValue $ws.Range("D21") "12.38"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("E22").Value = "  +0.31%  "

Set-TextValue $ws.Range("D23") "70.11"
$ws.Range("E23").Value = "  -0.04%  "

Set-TextValue $ws.Range("D24") "267.36"
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("E25").Value = "  -1.91%  "

Set-TextValue $ws.Range("D26") "8.32"
$ws.Range("E26").Value = "  +6.68%  "

Set-TextValue $ws.Range("D27") "27.00"
$ws.Range("E27").Value = "  +4.43%  "

Set-TextValue $ws.Range("D28") "0.171"
$ws.Range("E28").Value = "  +4.72%  "

Set-TextValue $ws.Range("D29") "7.21"
$ws.Range("E29").Value = "  -1.30%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  -2.11%  "

Set-TextValue $ws.Range("D32") "10.25"
$ws.Range("E32").Value = "  -0.15%  "

Set-TextValue $ws.Range("D33") "34.44"
$ws.Range("E33").Value = "  +0.41%  "

$ws.Range("E34").Value = "  +0.46%  "

Set-TextValue $ws.Range("D35") "50.49"
$ws.Range("E35").Value = "  -1.39%  "

$ws.Range("E36").Value = "  +2.24%  "

$ws.Range("E37").Value = "  -0.13%  "

Set-TextValue $ws.Range("D38") "3.34"
$ws.Range("E38").Value = "  +2.56%  "

$ws.Range("E39").Value = "  +8.13%  "

Set-TextValue $ws.Range("D40") "16.92"
$ws.Range("E40").Value = "  +2.81%  "

Set-TextValue $ws.Range("D41") "1.86"
$ws.Range("E41").Value = "  +2.18%  "

Set-TextValue $ws.Range("D42") "128.51"
$ws.Range("E42").Value = "  +3.20%  "

$ws.Range("E43").Value = "  -0.45%  "

Set-TextValue $ws.Range("D44") "2.54"
$ws.Range("E44").Value = "  +1.25%  "

Set-TextValue $ws.Range("D45") "3.70"
$ws.Range("E45").Value = "  +3.46%  "

Set-TextValue $ws.Range("D46") "21.81"
$ws.Range("E46").Value = "  +1.84%  "

Set-TextValue $ws.Range("D47") "2.47"
$ws.Range("E47").Value = "  +2.77%  "

$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").Value = "2.035.78"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("D50").Value = "3.348.26"
$ws.Range("E50").Value = "  +3.05%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.206"
$ws.Range("E51").Value = "  +7.13%  "
